$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.530.85"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.670.69"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'238.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").Value = "'0.06184"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").Value = "'0.07001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "1.668.72"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "'14.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.5913"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").Value = "'4.385"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "'75.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "'0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.9993"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "25.514.22"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "'0.000006773"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "1.880.85"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'8.741"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "'5.281"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'136.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("D26").Value = "'15.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'1.391"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'1.723"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.54%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'3.963"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.87%  "
$ws.Range("D31").Value = "'0.07812"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'3.662"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").Value = "'0.9985"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'0.04261"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "'2.601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "'0.6100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.28%  "
$ws.Range("D37").Value = "'0.9521"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'0.8590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").Value = "'0.9992"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'0.01487"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").Value = "'1.860"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").Value = "'95.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").Value = "'0.3777"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").Value = "'4.840"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'0.1119"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").Value = "'6.215"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Value = "'0.05249"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "'29.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").Value = "'1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.350"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
